$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: new columns I (I0) and J (IF), matching header style of H1 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-4: columns I and J both equal to 9 ---
$ws.Range("I2:J4").Value = 9
